$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 242, shifting existing rows 242:258 down to 243:259.
$ws.Rows.Item(242).Insert()

# Populate the new row 242 by copying formatting/values from the row now below it (243),
# then overwrite the cells that actually differ for the new record.
$ws.Range("A243:R243").Copy($ws.Range("A242:R242"))

$ws.Range("D242").Value = 45265
$ws.Range("J242").Value = 100
$ws.Range("K242").Value = 15000
$ws.Range("L242").Value = 16000
$ws.Range("M242").Value = 15500
$ws.Range("O242").Value = "Región de Arica y Parinacota"
$ws.Range("P242").Value = 258
